# "Added summary of work to PPT"
#
# Slide 19 ("Summary of Work", sldId 273) has a body textbox (shape id 175)
# that only contained a placeholder line break + "Blah blah". Replace it
# with the real four-part project summary, laid out as five paragraphs of
# text separated by blank paragraphs (mirroring how the author typed it in
# PowerPoint, with an Enter-Enter between each bullet-ish line).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(19)
$sh = $s.Shapes.Item("Google Shape;175;p30")
$tr = $sh.TextFrame.TextRange

$lines = @(
    "The creation of our webapp can be summarized into four parts:",
    "We acquired the data and stored it in SQL.",
    "We scripted functionality & design for the site in HTML and CSS.",
    "We enabled the site to pull from the APIs and our database using JavaScript.",
    "Finally, we mapped the locations & created the summary statistics in JavaScript."
)

$tr.Text = [string]::Join("`r`r", $lines) + "`r"
